$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.086.37'
$ws.Cells.Item(2, 5).Value = '  -1.03%  '

$ws.Cells.Item(3, 4).Value = '2.467.21'
$ws.Cells.Item(3, 5).Value = '  -1.22%  '

$ws.Cells.Item(4, 5).Value = '  +0.04%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '582.73'
$ws.Cells.Item(5, 5).Value = '  -1.61%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '166.87'
$ws.Cells.Item(6, 5).Value = '  -4.17%  '

$ws.Cells.Item(7, 5).Value = '  +0.09%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.513'
$ws.Cells.Item(8, 5).Value = '  -2.33%  '

$ws.Cells.Item(9, 4).Value = '2.468.28'
$ws.Cells.Item(9, 5).Value = '  -1.08%  '

$ws.Cells.Item(10, 5).Value = '  -4.41%  '

$ws.Cells.Item(11, 5).Value = '  -0.97%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '4.94'
$ws.Cells.Item(12, 5).Value = '  -3.17%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.332'
$ws.Cells.Item(13, 5).Value = '  -3.24%  '

$ws.Cells.Item(14, 4).Value = '2.928.47'
$ws.Cells.Item(14, 5).Value = '  -0.81%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '25.49'
$ws.Cells.Item(15, 5).Value = '  -3.23%  '

$ws.Cells.Item(16, 4).Value = '67.030.02'

$ws.Cells.Item(17, 5).Value = '  -4.79%  '

$ws.Cells.Item(18, 4).Value = '2.473.88'
$ws.Cells.Item(18, 5).Value = '  -1.23%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '11.31'
$ws.Cells.Item(19, 5).Value = '  -3.29%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '7.63'
$ws.Cells.Item(20, 5).Value = '  -4.49%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '355.52'
$ws.Cells.Item(21, 5).Value = '  -2.76%  '

$ws.Cells.Item(22, 5).Value = '  -2.61%  '

$ws.Cells.Item(23, 5).Value = '  +0.03%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '69.46'

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '4.23'
$ws.Cells.Item(25, 5).Value = '  -7.71%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '1.76'
$ws.Cells.Item(26, 5).Value = '  -8.16%  '

$ws.Cells.Item(27, 5).Value = '  -9.28%  '

$ws.Cells.Item(28, 5).Value = '  -0.39%  '

$ws.Cells.Item(29, 4).Value = '2.598.37'
$ws.Cells.Item(29, 5).Value = '  -0.78%  '

$ws.Cells.Item(30, 4).Value = '0.0₃0899'
$ws.Cells.Item(30, 5).Value = '  -6.83%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '512.07'
$ws.Cells.Item(31, 5).Value = '  -4.12%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '7.78'
$ws.Cells.Item(32, 5).Value = '  -6.37%  '

$ws.Cells.Item(33, 5).Value = '  -4.72%  '

$ws.Cells.Item(34, 5).Value = '  -6.15%  '

$ws.Cells.Item(35, 5).Value = '  +0.09%  '

$ws.Cells.Item(36, 5).Value = '  -6.72%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '158.51'
$ws.Cells.Item(37, 5).Value = '  -0.14%  '

$ws.Cells.Item(38, 5).Value = '  -0.30%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '18.44'
$ws.Cells.Item(39, 5).Value = '  -1.21%  '

$ws.Cells.Item(40, 5).Value = '  -6.38%  '

$ws.Cells.Item(41, 2).Value = 'USDe'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.00'
$ws.Cells.Item(41, 5).Value = '  +0.08%  '

$ws.Cells.Item(42, 2).Value = 'Stacks'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '1.66'
$ws.Cells.Item(42, 5).Value = '  -6.91%  '

$ws.Cells.Item(43, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.326'
$ws.Cells.Item(43, 5).Value = '  -6.90%  '

$ws.Cells.Item(44, 2).Value = 'RenderToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '4.77'
$ws.Cells.Item(44, 5).Value = '  -7.09%  '

$ws.Cells.Item(45, 2).Value = 'OKB'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '38.68'
$ws.Cells.Item(45, 5).Value = '  -2.88%  '

$ws.Cells.Item(46, 2).Value = 'dogwifhat'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '2.31'
$ws.Cells.Item(46, 5).Value = '  -8.27%  '

$ws.Cells.Item(47, 2).Value = 'Aave'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '141.24'
$ws.Cells.Item(47, 5).Value = '  -2.58%  '

$ws.Cells.Item(48, 2).Value = 'Filecoin'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '3.47'
$ws.Cells.Item(48, 5).Value = '  -5.80%  '

$ws.Cells.Item(49, 2).Value = 'ARBITRUM'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.514'
$ws.Cells.Item(49, 5).Value = '  -6.22%  '

$ws.Cells.Item(50, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(50, 4).Value = '0.0₆0253'
$ws.Cells.Item(50, 5).Value = '  -7.63%  '

$ws.Cells.Item(51, 2).Value = 'Optimism'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '1.59'
$ws.Cells.Item(51, 5).Value = '  -6.53%  '
